$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.229.70"
$ws.Range("E2").Value = "  -1.52%  "

# Row 3
$ws.Range("D3").Value = "1.851.52"
$ws.Range("E3").Value = "  -2.56%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.37"
$ws.Range("E5").Value = "  -2.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4693"
$ws.Range("E7").Value = "  -2.58%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2716"
$ws.Range("E8").Value = "  -4.50%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06380"
$ws.Range("E9").Value = "  -2.83%  "

# Row 10
$ws.Range("D10").Value = "1.854.23"
$ws.Range("E10").Value = "  -3.60%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07421"
$ws.Range("E11").Value = "  -0.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.27"
$ws.Range("E12").Value = "  -2.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.940"
$ws.Range("E13").Value = "  -3.56%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.09"
$ws.Range("E14").Value = "  -3.55%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6276"
$ws.Range("E15").Value = "  -6.01%  "

# Row 16
$ws.Range("D16").Value = "30.178.70"
$ws.Range("E16").Value = "  -1.60%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.71"
$ws.Range("E18").Value = "  -2.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.59"
$ws.Range("E19").Value = "  -5.57%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007328"
$ws.Range("E20").Value = "  -3.85%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.934"
$ws.Range("E22").Value = "  -6.99%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.932"
$ws.Range("E23").Value = "  -4.78%  "

# Row 24
$ws.Range("E24").Value = "  -0.81%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.53"
$ws.Range("E25").Value = "  -2.22%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.80"
$ws.Range("E26").Value = "  -4.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.870"
$ws.Range("E27").Value = "  -4.80%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1034"
$ws.Range("E28").Value = "  +2.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.381"
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.118"
$ws.Range("E30").Value = "  -5.86%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.874"
$ws.Range("E31").Value = "  -4.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04893"
$ws.Range("E32").Value = "  -3.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.154"
$ws.Range("E33").Value = "  -5.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7108"
$ws.Range("E34").Value = "  -5.86%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9997"
$ws.Range("E35").Value = "  -0.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.696"
$ws.Range("E36").Value = "  -0.70%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01856"
$ws.Range("E37").Value = "  -1.55%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.639"
$ws.Range("E38").Value = "  -0.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9052"
$ws.Range("E39").Value = "  -1.42%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.948"
$ws.Range("E40").Value = "  -6.60%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.10"
$ws.Range("E41").Value = "  -1.97%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9981"
$ws.Range("E42").Value = "  -0.69%  "

# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.551"
$ws.Range("E43").Value = "  -3.94%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4083"
$ws.Range("E44").Value = "  -5.29%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.055"
$ws.Range("E45").Value = "  -5.15%  "

# Row 46
$ws.Range("E46").Value = "  -7.10%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1189"
$ws.Range("E47").Value = "  -6.85%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.600"
$ws.Range("E48").Value = "  -4.06%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.10"
$ws.Range("E49").Value = "  -2.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.387"
$ws.Range("E50").Value = "  -6.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05566"

